# Fix tests for list codec: correct a few bit-pattern truth-table rows,
# collapse ipv4/ipv6/geo rows into a single "open" row, drop the now-unused
# [long]/[float] example labels on rows 17/18, and relabel row 16's example
# from "integer" to "[float]".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 (datetime) : flip the C/D bits ---
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1

# --- Row 16 : was "ipv4", becomes "geo"; flip B/E/F bits; relabel example ---
$ws.Range("A16").Value = "geo"
$ws.Range("B16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("K16").Value = "[float]"

# --- Row 17 : was "ipv6", becomes "open"; clear the B:F bit flags and the example ---
$ws.Range("A17").Value = "open"
$ws.Range("B17:F17").ClearContents()
$ws.Range("K17").ClearContents()

# --- Row 18 : was "geo", becomes "open"; clear the B:F bit flags and the example ---
$ws.Range("A18").Value = "open"
$ws.Range("B18:F18").ClearContents()
$ws.Range("K18").ClearContents()

# --- View state: zoom to 130% and move the selection to K18 ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("K18").Select()
